$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in column A (same as the sheet's data extent)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Columns that carry the HYPERLINK() formulas which need a friendly
# display-text second argument equal to the row's "Beteckning" (column A).
$linkCols = @(19, 20, 22, 23, 24, 25)  # S, T, V, W, X, Y

for ($r = 2; $r -le $lastRow; $r++) {

    # Update the "Förändrad" (changed) date in column C for every data row.
    $cCell = $ws.Cells.Item($r, 3)
    if ($cCell.Value2 -ne $null) {
        $cCell.Value = 45186
    }

    # Designation text used as the HYPERLINK friendly name.
    $beteckning = $ws.Cells.Item($r, 1).Value2

    foreach ($col in $linkCols) {
        $cell = $ws.Cells.Item($r, $col)
        $formula = $cell.Formula
        if ([string]::IsNullOrEmpty($formula)) {
            continue
        }
        if ($formula.StartsWith("=HYPERLINK(") -and -not $formula.Contains(",")) {
            $newFormula = $formula.Substring(0, $formula.Length - 1) + ', "' + $beteckning + '")'
            $cell.Formula = $newFormula
        }
    }
}
